$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'40.110.72"
$ws.Range("E2").Value = "  +0.84%  "

$ws.Range("D3").Value = "'2.239.84"
$ws.Range("E3").Value = "  -3.45%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'294.98"
$ws.Range("E5").Value = "  -4.42%  "

$ws.Range("D6").Value = "'86.91"
$ws.Range("E6").Value = "  +4.87%  "

$ws.Range("D7").Value = "'0.514"
$ws.Range("E7").Value = "  -1.45%  "

$ws.Range("D9").Value = "'0.472"
$ws.Range("E9").Value = "  -0.58%  "

$ws.Range("B10").Value = "'Avalanche"
$ws.Range("C10").Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D10").Value = "'30.81"
$ws.Range("E10").Value = "  +5.06%  "

$ws.Range("B11").Value = "'Dogecoin"
$ws.Range("C11").Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.0799"
$ws.Range("E11").Value = "  +0.07%  "

$ws.Range("D12").Value = "'47.29"

$ws.Range("E13").Value = "  -1.87%  "

$ws.Range("D14").Value = "'6.41"
$ws.Range("E14").Value = "  +1.59%  "

$ws.Range("D15").Value = "'2.581.07"
$ws.Range("E15").Value = "  -3.72%  "

$ws.Range("D16").Value = "'14.23"
$ws.Range("E16").Value = "  -2.20%  "

$ws.Range("D17").Value = "'2.233.79"
$ws.Range("E17").Value = "  -4.16%  "

$ws.Range("D18").Value = "'0.728"
$ws.Range("E18").Value = "  -2.55%  "

$ws.Range("D19").Value = "'39.977.52"
$ws.Range("E19").Value = "  +0.66%  "

$ws.Range("D20").Value = "'0.0₃0895"
$ws.Range("E20").Value = "  +0.60%  "

$ws.Range("D21").Value = "'5.81"
$ws.Range("E21").Value = "  -3.03%  "

$ws.Range("D22").Value = "'10.78"
$ws.Range("E22").Value = "  +3.95%  "

$ws.Range("D23").Value = "'65.60"
$ws.Range("E23").Value = "  -3.13%  "

$ws.Range("D24").Value = "'235.05"
$ws.Range("E24").Value = "  +0.64%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D26").Value = "'2.44"
$ws.Range("E26").Value = "  -2.77%  "

$ws.Range("D27").Value = "'1.85"
$ws.Range("E27").Value = "  +3.30%  "

$ws.Range("D28").Value = "'23.09"
$ws.Range("E28").Value = "  -0.32%  "

$ws.Range("E29").Value = "  +0.96%  "

$ws.Range("D30").Value = "'9.27"
$ws.Range("E30").Value = "  +1.57%  "

$ws.Range("D31").Value = "'33.60"
$ws.Range("E31").Value = "  +0.15%  "

$ws.Range("D32").Value = "'154.98"
$ws.Range("E32").Value = "  +1.94%  "

$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.29%  "

$ws.Range("D34").Value = "'4.88"
$ws.Range("E34").Value = "  -2.72%  "

$ws.Range("D35").Value = "'0.0715"
$ws.Range("E35").Value = "  +1.39%  "

$ws.Range("E36").Value = "  -2.88%  "

$ws.Range("D37").Value = "'16.67"
$ws.Range("E37").Value = "  +9.83%  "

$ws.Range("E38").Value = "  -0.49%  "

$ws.Range("D39").Value = "'0.0995"
$ws.Range("E39").Value = "  +2.44%  "

$ws.Range("D40").Value = "'2.70"
$ws.Range("E40").Value = "  -0.88%  "

$ws.Range("D41").Value = "'1.68"
$ws.Range("E41").Value = "  +0.72%  "

$ws.Range("D42").Value = "'3.80"
$ws.Range("E42").Value = "  +2.37%  "

$ws.Range("D43").Value = "'1.955.44"
$ws.Range("E43").Value = "  -0.44%  "

$ws.Range("D44").Value = "'2.19"
$ws.Range("E44").Value = "  -3.02%  "

$ws.Range("E45").Value = "  +4.27%  "

$ws.Range("D46").Value = "'9.54"
$ws.Range("E46").Value = "  +2.18%  "

$ws.Range("E47").Value = "  -4.88%  "

$ws.Range("D48").Value = "'2.62"
$ws.Range("E48").Value = "  -0.26%  "

$ws.Range("D49").Value = "'2.452.37"
$ws.Range("E49").Value = "  -3.53%  "

$ws.Range("D50").Value = "'71.23"
$ws.Range("E50").Value = "  +2.87%  "

$ws.Range("E51").Value = "  +9.76%  "
